# Auto-generated edit script: updates leve-profit figures per the commit diff.
# Values were recomputed by the scheduled pricing runner; this script writes
# the refreshed numbers (and clears now-empty profit cells) back into each sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 46.037037
$ws.Range("I2").Value = 36.307693
$ws.Range("K2").Value = 36.307693
$ws.Range("M2").Value = 76.692307

$ws.Range("H6").Value = 2463200.8
$ws.Range("I6").Value = 2463200.8
$ws.Range("K6").Value = 7389602.399999999
$ws.Range("M6").Value = -7389490.399999999

$ws.Range("H15").Value = 48804
$ws.Range("I15").Value = 48804
$ws.Range("K15").Value = 146412
$ws.Range("M15").Value = -146243

$ws.Range("H86").Value = 100002820
$ws.Range("I86").Value = 131315064
$ws.Range("J86").Value = 13894139
$ws.Range("K86").Value = 131315064
$ws.Range("L86").Value = 13894139
$ws.Range("M86").Value = -131313941
$ws.Range("N86").Value = -13896385

$ws.Range("H89").Value = 100002820
$ws.Range("I89").Value = 131315064
$ws.Range("J89").Value = 13894139
$ws.Range("K89").Value = 656575320
$ws.Range("L89").Value = 69470695
$ws.Range("M89").Value = -656569704
$ws.Range("N89").Value = -69481927

$ws.Range("H113").Value = 61121710
$ws.Range("I113").Value = 55557060
$ws.Range("J113").Value = 62512876
$ws.Range("K113").Value = 55557060
$ws.Range("L113").Value = 62512876
$ws.Range("M113").Value = -55553806
$ws.Range("N113").Value = -62519384

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()

$ws.Range("H45").Value = 2823.3809
$ws.Range("J45").Value = 6688
$ws.Range("L45").Value = 6688
$ws.Range("N45").Value = -7442

$ws.Range("H61").Value = 4718.9
$ws.Range("I61").Value = 2266.9167
$ws.Range("K61").Value = 2266.9167
$ws.Range("M61").Value = -2054.9167

$ws.Range("H76").Value = 41122
$ws.Range("J76").Value = 41122
$ws.Range("L76").Value = 41122
$ws.Range("N76").Value = -41798

$ws.Range("H79").Value = 41122
$ws.Range("J79").Value = 41122
$ws.Range("L79").Value = 41122
$ws.Range("N79").Value = -43462

$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").ClearContents()

$ws.Range("H132").Value = 7594.72
$ws.Range("I132").Value = 5903.091
$ws.Range("J132").Value = 20000
$ws.Range("K132").Value = 17709.273
$ws.Range("L132").Value = 60000
$ws.Range("M132").Value = -15179.273
$ws.Range("N132").Value = -65060

$ws.Range("H136").Value = 4718.9
$ws.Range("I136").Value = 2266.9167
$ws.Range("K136").Value = 6800.750100000001
$ws.Range("M136").Value = -4250.750100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()

$ws.Range("H50").Value = 52178
$ws.Range("J50").Value = 52178
$ws.Range("L50").Value = 52178
$ws.Range("N50").Value = -53326

$ws.Range("H86").Value = 65845700
$ws.Range("I86").Value = 25102520
$ws.Range("J86").Value = 111115896
$ws.Range("K86").Value = 25102520
$ws.Range("L86").Value = 111115896
$ws.Range("M86").Value = -25101397
$ws.Range("N86").Value = -111118142

$ws.Range("H89").Value = 65845700
$ws.Range("I89").Value = 25102520
$ws.Range("J89").Value = 111115896
$ws.Range("K89").Value = 125512600
$ws.Range("L89").Value = 555579480
$ws.Range("M89").Value = -125506984
$ws.Range("N89").Value = -555590712

$ws.Range("H130").Value = 87000
$ws.Range("J130").Value = 87000
$ws.Range("L130").Value = 87000
$ws.Range("N130").Value = -97040

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5909.85
$ws.Range("I16").Value = 3601.1428
$ws.Range("J16").Value = 7153
$ws.Range("K16").Value = 3601.1428
$ws.Range("L16").Value = 7153
$ws.Range("M16").Value = -3314.1428
$ws.Range("N16").Value = -7727

$ws.Range("H31").Value = 9741.200000000001
$ws.Range("I31").Value = 4150.385
$ws.Range("J31").Value = 12433.074
$ws.Range("K31").Value = 4150.385
$ws.Range("L31").Value = 12433.074
$ws.Range("M31").Value = -3855.385
$ws.Range("N31").Value = -13023.074

$ws.Range("H34").Value = 9741.200000000001
$ws.Range("I34").Value = 4150.385
$ws.Range("J34").Value = 12433.074
$ws.Range("K34").Value = 4150.385
$ws.Range("L34").Value = 12433.074
$ws.Range("M34").Value = -3948.385
$ws.Range("N34").Value = -12837.074

$ws.Range("H113").Value = 5909.85
$ws.Range("I113").Value = 3601.1428
$ws.Range("J113").Value = 7153
$ws.Range("K113").Value = 3601.1428
$ws.Range("L113").Value = 7153
$ws.Range("M113").Value = -1431.1428
$ws.Range("N113").Value = -11493

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3337162.2
$ws.Range("I5").Value = 5715792.5
$ws.Range("K5").Value = 17147377.5
$ws.Range("M5").Value = -17147265.5

$ws.Range("H129").Value = 18520114
$ws.Range("J129").Value = 27779352
$ws.Range("L129").Value = 83338056
$ws.Range("N129").Value = -83348056

$ws.Range("H131").Value = 1884.1892
$ws.Range("I131").Value = 1264.5454
$ws.Range("J131").Value = 2146.3462
$ws.Range("K131").Value = 3793.6362
$ws.Range("L131").Value = 6439.0386
$ws.Range("M131").Value = 1246.3638
$ws.Range("N131").Value = -16519.0386

$ws.Range("H135").Value = 3337162.2
$ws.Range("I135").Value = 5715792.5
$ws.Range("K135").Value = 51442132.5
$ws.Range("M135").Value = -51439597.5

$ws.Range("H139").Value = 5294.885
$ws.Range("I139").Value = 2083.5
$ws.Range("J139").Value = 9041.5
$ws.Range("K139").Value = 6250.5
$ws.Range("L139").Value = 27124.5
$ws.Range("M139").Value = -1110.5
$ws.Range("N139").Value = -37404.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H45").Value = 24000
$ws.Range("J45").Value = 24000
$ws.Range("L45").Value = 24000
$ws.Range("N45").Value = -25118

$ws.Range("H130").Value = 88983
$ws.Range("J130").Value = 88983
$ws.Range("L130").Value = 88983
$ws.Range("N130").Value = -99023

$ws.Range("H132").Value = 7699.8
$ws.Range("J132").Value = 5750
$ws.Range("L132").Value = 17250
$ws.Range("N132").Value = -22310

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 15161481
$ws.Range("I132").Value = 41674708
$ws.Range("J132").Value = 11066.619
$ws.Range("K132").Value = 125024124
$ws.Range("L132").Value = 33199.857
$ws.Range("M132").Value = -125021594
$ws.Range("N132").Value = -38259.857

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 68479.5
$ws.Range("J16").Value = 68479.5
$ws.Range("L16").Value = 68479.5
$ws.Range("N16").Value = -69063.5

$ws.Range("H126").Value = 5142
$ws.Range("I126").Value = 2998.75
$ws.Range("K126").Value = 8996.25
$ws.Range("M126").Value = -6526.25

$ws.Range("H136").Value = 30337830
$ws.Range("I136").Value = 111111920
$ws.Range("K136").Value = 333335760
$ws.Range("M136").Value = -333333210
